$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.943.66'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '2.041.61'
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.32'
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.43'
$ws.Range("E8").Value = '  -1.11%  '

$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0780'
$ws.Range("E10").Value = '  +2.52%  '

$ws.Range("E11").Value = '  +1.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.67'
$ws.Range("E12").Value = '  +3.33%  '

$ws.Range("D13").Value = '2.339.92'
$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.56'
$ws.Range("E14").Value = '  +5.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.786'
$ws.Range("E15").Value = '  -5.11%  '

$ws.Range("D16").Value = '2.043.21'
$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").Value = '36.874.12'
$ws.Range("E17").Value = '  -0.46%  '

$ws.Range("E18").Value = '  +13.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.41'
$ws.Range("E19").Value = '  +1.16%  '

$ws.Range("D20").Value = '0.0₃0889'
$ws.Range("E20").Value = '  +2.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.29'
$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '234.98'
$ws.Range("E22").Value = '  -1.52%  '

$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.34'
$ws.Range("E24").Value = '  -3.14%  '

$ws.Range("E25").Value = '  +7.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.08'
$ws.Range("E26").Value = '  -1.90%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("E27").Value = '  -1.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.64'
$ws.Range("E28").Value = '  -3.54%  '

$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("E30").Value = '  +1.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.64'
$ws.Range("E31").Value = '  +1.46%  '

$ws.Range("E32").Value = '  -3.94%  '

$ws.Range("E33").Value = '  -0.26%  '

$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("E35").Value = '  +0.62%  '

$ws.Range("E36").Value = '  -4.53%  '

$ws.Range("E37").Value = '  -1.23%  '

$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("E38").Value = '  +15.51%  '

$ws.Range("E39").Value = '  -2.12%  '

$ws.Range("B40").Value = 'Cronos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.105'
$ws.Range("E40").Value = '  -2.27%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0219'
$ws.Range("E41").Value = '  -3.10%  '

$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.85'
$ws.Range("E42").Value = '  +21.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.13'
$ws.Range("E43").Value = '  -6.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.86'
$ws.Range("E44").Value = '  -3.29%  '

$ws.Range("E45").Value = '  -3.64%  '

$ws.Range("E46").Value = '  +0.45%  '

$ws.Range("D47").Value = '1.270.24'
$ws.Range("E47").Value = '  -2.76%  '

$ws.Range("E48").Value = '  -2.05%  '

$ws.Range("D49").Value = '2.222.03'
$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.64'
$ws.Range("E50").Value = '  -3.95%  '

$ws.Range("E51").Value = '  -7.56%  '
